# Weekly refresh: re-sort / reshuffle the daily price rows (rows 3-19) of the
# "Hortaliza, Vega Monumental Concepción - Espárragos" sheet.
#
# Every data column except A,B,C,E,F,G,Q,R (which are constant across all
# rows for this sheet) moves from its old row to a new row according to the
# mapping below: new row N receives the data that used to live in old row
# mapping[N].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> old row
$mapping = @{
    3  = 5
    4  = 10
    5  = 11
    6  = 17
    7  = 6
    8  = 15
    9  = 9
    10 = 8
    11 = 7
    12 = 3
    13 = 19
    14 = 12
    15 = 16
    16 = 4
    17 = 13
    18 = 14
    19 = 18
}

# Columns that vary row to row: D, H, I, J, K, L, M, N, O, P
$cols = @(4, 8, 9, 10, 11, 12, 13, 14, 15, 16)

# Read every old row's values into memory first (this is a full permutation,
# so we must snapshot everything before writing anything).
$snapshot = @{}
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    if (-not $snapshot.ContainsKey($oldRow)) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($oldRow, $c).Value2
        }
        $snapshot[$oldRow] = $rowVals
    }
}

# Now write the snapshots into their new row positions.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $rowVals = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value = $rowVals[$c]
    }
}

Write-Output "done"
